$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Links"
$ws.Range("D2").Value = "https://film-oneri.streamlit.app/"

$ws.Columns.Item(4).ColumnWidth = 27.88671875

$ws.Range("B11").Select()
